$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple pairwise row-content swaps (columns B..AC); column A (row index) stays fixed.
$rowA = $ws.Range("B23:AC23")
$rowB = $ws.Range("B24:AC24")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B38:AC38")
$rowB = $ws.Range("B39:AC39")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B54:AC54")
$rowB = $ws.Range("B55:AC55")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B72:AC72")
$rowB = $ws.Range("B73:AC73")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B74:AC74")
$rowB = $ws.Range("B75:AC75")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B82:AC82")
$rowB = $ws.Range("B83:AC83")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B96:AC96")
$rowB = $ws.Range("B97:AC97")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B107:AC107")
$rowB = $ws.Range("B108:AC108")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B109:AC109")
$rowB = $ws.Range("B110:AC110")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B114:AC114")
$rowB = $ws.Range("B115:AC115")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B131:AC131")
$rowB = $ws.Range("B132:AC132")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B137:AC137")
$rowB = $ws.Range("B138:AC138")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B142:AC142")
$rowB = $ws.Range("B143:AC143")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B148:AC148")
$rowB = $ws.Range("B149:AC149")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B169:AC169")
$rowB = $ws.Range("B170:AC170")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B172:AC172")
$rowB = $ws.Range("B174:AC174")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B185:AC185")
$rowB = $ws.Range("B186:AC186")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

$rowA = $ws.Range("B195:AC195")
$rowB = $ws.Range("B196:AC196")
$valA = $rowA.Value2
$valB = $rowB.Value2
$rowA.Value2 = $valB
$rowB.Value2 = $valA

# 3-way rotation among rows 31, 33, 34 (row32 untouched):
#   row31 <- original row33, row33 <- original row34, row34 <- original row31
$v31 = $ws.Range("B31:AC31").Value2
$v33 = $ws.Range("B33:AC33").Value2
$v34 = $ws.Range("B34:AC34").Value2
$ws.Range("B31:AC31").Value2 = $v33
$ws.Range("B33:AC33").Value2 = $v34
$ws.Range("B34:AC34").Value2 = $v31
